# "Fixed typos ... Kraken2 -> Kraken 2"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header/source row (row 2) held the literal text "Kraken2" - fix the typo.
$ws.Range("Z2").Value = "Kraken 2"

# Rows 3-13 held independent copies of that same text; turn them into
# formulas that reference the (now corrected) source cell instead, so they
# can never drift out of sync again.
for ($r = 3; $r -le 13; $r++) {
    $ws.Cells.Item($r, 26).Formula = "=Z`$2"
}

# Move the active selection to AA15, matching the author's final cursor
# position when they saved the file.
$ws.Range("AA15").Select()
